# Generate Report for Handback
#
# - Overview sheet: zh-cn / de-de status columns move from
#   "Ready for handoff" to "Handed back: in sync with en-US".
# - zh-cn / de-de sheets: the "Latest Target File" (J) and
#   "Latest Handback File" (K) columns get populated (J as a hyperlink to
#   the source .md doc, K with the generated xliff file name), and the
#   "Latest Handback DateTime" (L) timestamp is stamped.
# - Column widths are widened on the affected columns to fit the new text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22c6efe221dbb63e8e304370850661200f83e506/e2e/3c76a564-2bc0-40fa-bdd9-7f672272ffa3.md", "", "", "3c76a564-2bc0-40fa-bdd9-7f672272ffa3.md")
$zhcn.Range("K2").Value = "3c76a564-2bc0-40fa-bdd9-7f672272ffa3.cf39bda889242f13e1997b684b7d3605a8e2d062.zh-cn.xlf"
$zhcn.Range("L2").Value = "2017-02-17 10:22:13"

$zhcn.Hyperlinks.Add($zhcn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22c6efe221dbb63e8e304370850661200f83e506/e2e/56a178da-5ee8-40ec-b336-48b041c3ca5e.md", "", "", "56a178da-5ee8-40ec-b336-48b041c3ca5e.md")
$zhcn.Range("K3").Value = "56a178da-5ee8-40ec-b336-48b041c3ca5e.66aeac91f46b606504cec4445255b6445d0faaa6.zh-cn.xlf"
$zhcn.Range("L3").Value = "2017-02-17 10:22:13"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15
$zhcn.Columns.Item(11).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22c6efe221dbb63e8e304370850661200f83e506/e2e/3c76a564-2bc0-40fa-bdd9-7f672272ffa3.md", "", "", "3c76a564-2bc0-40fa-bdd9-7f672272ffa3.md")
$dede.Range("K2").Value = "3c76a564-2bc0-40fa-bdd9-7f672272ffa3.cf39bda889242f13e1997b684b7d3605a8e2d062.de-de.xlf"
$dede.Range("L2").Value = "2017-02-17 10:22:35"

$dede.Hyperlinks.Add($dede.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22c6efe221dbb63e8e304370850661200f83e506/e2e/56a178da-5ee8-40ec-b336-48b041c3ca5e.md", "", "", "56a178da-5ee8-40ec-b336-48b041c3ca5e.md")
$dede.Range("K3").Value = "56a178da-5ee8-40ec-b336-48b041c3ca5e.66aeac91f46b606504cec4445255b6445d0faaa6.de-de.xlf"
$dede.Range("L3").Value = "2017-02-17 10:22:35"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(10).ColumnWidth = 39.15
$dede.Columns.Item(11).ColumnWidth = 39.15
